$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C5").Value = "mfrd"
$ws.Range("C6").Value = "afrd"

$ws.Range("D5").Value = 1
$ws.Range("E5").Value = 0

$ws.Range("D6").Value = 0.54231
$ws.Range("E6").Value = 0.05329
